# Daily attendance processing - reorder "Recorded By" (column G) names so
# that any "System" / "system" token(s) are moved to the front of the
# comma-separated list, preserving the relative order of all other
# entries (exact-case "System" first, then other case variants of
# "system", then the remaining names in their original order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-RecordedBy {
    param([string]$val)

    if ([string]::IsNullOrEmpty($val)) {
        return $val
    }

    $rawParts = $val -split ","
    $trimmed = @()
    foreach ($p in $rawParts) {
        $trimmed += $p.Trim()
    }

    $hasSystem = $false
    foreach ($t in $trimmed) {
        if ($t.ToLower().Equals("system")) {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) {
        return $val
    }

    $exact = @()
    $ciOnly = @()
    $others = @()
    foreach ($t in $trimmed) {
        if ($t.Equals("System")) {
            $exact += $t
        } elseif ($t.ToLower().Equals("system")) {
            $ciOnly += $t
        } else {
            $others += $t
        }
    }

    $ordered = $exact + $ciOnly + $others
    return [string]::Join(", ", $ordered)
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2
    if ($current -ne $null) {
        $updated = Reorder-RecordedBy $current
        if (-not $updated.Equals($current)) {
            $cell.Value2 = $updated
        }
    }
}
